$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepend the NACE letter code (column A) to the sector description (column B)
# for each data row, e.g. "Agriculture, Forestry, and Fishing" -> "A. Agriculture, Forestry, and Fishing".
for ($row = 2; $row -le 23; $row++) {
    $code = $ws.Cells.Item($row, 1).Value()
    $desc = $ws.Cells.Item($row, 2).Value()
    $ws.Cells.Item($row, 2).Value = ($code + ". " + $desc)
}

# Scroll the sheet view so row 7 is the top-left visible row instead of row 13.
$excel.ActiveWindow.ScrollRow = 7
